$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Median Value (C) and Tier (D) columns for rows 2-80
# based on the updated poverty table calculation methodology.
$values = @{
    2 = @(0.8652770716827656, "Below Median")
    3 = @(0.6507371631926793, "Below Median")
    4 = @(2.153736654804271, "1st Tier")
    5 = @(1.544483985765124, "2nd Tier")
    6 = @(2.109134045077105, "1st Tier")
    7 = @(1.525165226232842, "2nd Tier")
    8 = @(1.020640569395018, "4th Tier")
    9 = @(2.763938315539739, "1st Tier")
    10 = @(0.6761565836298933, "Below Median")
    11 = @(0.8640569395017794, "Below Median")
    12 = @(0.9150991357397052, "Below Median")
    13 = @(1, "4th Tier")
    14 = @(1.138790035587189, "4th Tier")
    15 = @(1.175563463819691, "3rd Tier")
    16 = @(0.896797153024911, "Below Median")
    17 = @(0.9306049822064056, "Below Median")
    18 = @(0.8125741399762753, "Below Median")
    19 = @(0.8390668248319494, "Below Median")
    20 = @(0.697508896797153, "Below Median")
    21 = @(0.6790035587188612, "Below Median")
    22 = @(0.6456532791052364, "Below Median")
    23 = @(0.6751398068124047, "Below Median")
    24 = @(1.928825622775801, "1st Tier")
    25 = @(1.653279105236401, "2nd Tier")
    26 = @(1.962040332147094, "1st Tier")
    27 = @(1.421708185053381, "2nd Tier")
    28 = @(1.405693950177936, "3rd Tier")
    29 = @(0.706049822064057, "Below Median")
    30 = @(1.6355871886121, "2nd Tier")
    31 = @(1.104784499802294, "4th Tier")
    32 = @(1.829181494661922, "1st Tier")
    33 = @(0.7995255041518387, "Below Median")
    34 = @(1.369598373157092, "3rd Tier")
    35 = @(0.5599051008303678, "Below Median")
    36 = @(1.638027452974072, "2nd Tier")
    37 = @(0.6405693950177936, "Below Median")
    38 = @(1.275038129130656, "3rd Tier")
    39 = @(1.206914082358922, "3rd Tier")
    40 = @(1.752313167259786, "1st Tier")
    41 = @(1.019572953736655, "4th Tier")
    42 = @(0.9395017793594306, "Below Median")
    43 = @(1.579181494661922, "2nd Tier")
    44 = @(1.076027175671304, "4th Tier")
    45 = @(1.049110320284698, "4th Tier")
    46 = @(0.4427046263345196, "Below Median")
    47 = @(0.8315539739027283, "Below Median")
    48 = @(1.164590747330961, "3rd Tier")
    49 = @(1.060498220640569, "4th Tier")
    50 = @(0.9635231316725978, "Below Median")
    51 = @(0.604982206405694, "Below Median")
    52 = @(0.7288256227758008, "Below Median")
    53 = @(1.67497034400949, "1st Tier")
    54 = @(0.599644128113879, "Below Median")
    55 = @(1.234367056431113, "3rd Tier")
    56 = @(0.8784951703101169, "Below Median")
    57 = @(0.8042704626334519, "Below Median")
    58 = @(0.5701067615658363, "Below Median")
    59 = @(0.2740213523131673, "Below Median")
    60 = @(0.523640061006609, "Below Median")
    61 = @(1.238434163701068, "3rd Tier")
    62 = @(1.921708185053381, "1st Tier")
    63 = @(0.6818505338078291, "Below Median")
    64 = @(0.6220640569395018, "Below Median")
    65 = @(0.6423487544483986, "Below Median")
    66 = @(0.6749703440094899, "Below Median")
    67 = @(1.268327402135231, "3rd Tier")
    68 = @(1.421708185053381, "2nd Tier")
    69 = @(0.4498220640569395, "Below Median")
    70 = @(0.8042704626334519, "Below Median")
    71 = @(1.134519572953737, "4th Tier")
    72 = @(1.412099644128114, "2nd Tier")
    73 = @(0.604982206405694, "Below Median")
    74 = @(1.218098627351296, "3rd Tier")
    75 = @(2.542348754448398, "1st Tier")
    76 = @(1.111506524317912, "4th Tier")
    77 = @(0.7871886120996441, "Below Median")
    78 = @(1.50711743772242, "2nd Tier")
    79 = @(0.603202846975089, "Below Median")
    80 = @(0.8291814946619217, "Below Median")
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 4).Value = $pair[1]
}
